# "Generate Report for Handback"
#
# The handback transform failed for file 0c8ae761-4596-4822-9a7c-ec2bec2c275a
# (its handback archive name didn't match the handoff file name), so the
# status text changes from "Ready for handoff" to "Handback transform
# failed" everywhere that row's status is shown, and the per-language
# "Error Detail" column gets populated with the mismatch description.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$failedStatus = "Handback transform failed"

# Row 3 everywhere is the 0c8ae761-... file; update its status text.
$overview.Range("E3").Value = $failedStatus
$overview.Range("F3").Value = $failedStatus
$zhcn.Range("C3").Value = $failedStatus
$dede.Range("C3").Value = $failedStatus

# Populate the (previously empty) "Error Detail" column for that row on
# each language sheet, and widen the column so the message is readable.
# (ColumnWidth has a fixed +0.8333... offset baked into the stored OOXML
# width vs. the character-count value Excel shows in the UI, so back that
# out here to land on an exact stored width of 40.)
$targetStoredWidth = 40
$columnWidthOffset = 0.8333333333333334

$zhcn.Range("P3").Value = "Handback file name: 01aw42g4.vd5 is different with handoff file name: 0c8ae761-4596-4822-9a7c-ec2bec2c275a.e10d9dfcae65695686087f1da84cc0448e18879d.zh-cn."
$zhcn.Columns.Item(16).ColumnWidth = $targetStoredWidth - $columnWidthOffset

$dede.Range("P3").Value = "Handback file name: 01aw42g4.vd5 is different with handoff file name: 0c8ae761-4596-4822-9a7c-ec2bec2c275a.e10d9dfcae65695686087f1da84cc0448e18879d.de-de."
$dede.Columns.Item(16).ColumnWidth = $targetStoredWidth - $columnWidthOffset
